$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-03 10:55:00"
$wsZhCn.Range("G14").Value = "2016-03-03 10:56:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-03 10:55:14"
$wsDeDe.Range("G14").Value = "2016-03-03 10:56:24"
